$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 BaggingClassifier(estimator=SVC(C=0.25,
                                                 class_weight='balanced',
                                                 random_state=42),
                                   n_estimators=100, random_state=42))])"
$ws.Range("B2").Value = 0.5298740148740149
$ws.Range("C2").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': MinMaxScaler(), 'model__n_estimators': 100, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 0.25}"
$ws.Range("D2").Value = 0.6540104805010939
$ws.Range("E2").Value = 0.4704425666925667
$ws.Range("F2").Value = 0.7333333333333334
$ws.Range("G2").Value = 0.8270104155185082
$ws.Range("H2").Value = 0.590218253968254
$ws.Range("I2").Value = 0.7857142857142857
$ws.Range("J2").Value = 0.5703900709219858
$ws.Range("K2").Value = 0.4430555555555555
$ws.Range("L2").Value = 0.6875
$ws.Range("M2").Value = "[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]"
$ws.Range("N2").Value = "[1 0 1 1 1 1 1 0 1 0 0 1 0 1 0 0 0 1 0 1 1 1 1 0]"
$ws.Range("O2").Value = 42

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7fa6283035e0>),
                ('model',
                 BaggingClassifier(estimator=SVC(C=0.25,
                                                 class_weight='balanced',
                                                 random_state=42),
                                   n_estimators=100, random_state=42))])"
$ws.Range("B3").Value = 0.5494444444444444
$ws.Range("C3").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fa6ccf7cb50>, 'scaler': MinMaxScaler(), 'model__n_estimators': 100, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 0.25}"
$ws.Range("D3").Value = 0.6738050314870079
$ws.Range("E3").Value = 0.4856897731897732
$ws.Range("F3").Value = 0.6896551724137931
$ws.Range("G3").Value = 0.80233932078281
$ws.Range("H3").Value = 0.6626256613756614
$ws.Range("I3").Value = 0.7692307692307693
$ws.Range("J3").Value = 0.5978723404255321
$ws.Range("K3").Value = 0.4166666666666667
$ws.Range("L3").Value = 0.625
$ws.Range("M3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]"
$ws.Range("N3").Value = "[0 1 0 1 1 0 0 0 1 0 0 0 1 0 1 1 1 1 1 0 1 0 1 1]"
$ws.Range("O3").Value = 69

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 BaggingClassifier(estimator=SVC(C=0.25,
                                                 class_weight='balanced',
                                                 random_state=42),
                                   n_estimators=50, random_state=42))])"
$ws.Range("B4").Value = 0.5890909090909091
$ws.Range("C4").Value = "{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': MinMaxScaler(), 'model__n_estimators': 50, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 0.25}"
$ws.Range("D4").Value = 0.6541430460895407
$ws.Range("E4").Value = 0.5333589096089096
$ws.Range("F4").Value = 0.6451612903225806
$ws.Range("G4").Value = 0.7962490173887399
$ws.Range("H4").Value = 0.649593253968254
$ws.Range("I4").Value = 0.8333333333333334
$ws.Range("J4").Value = 0.5648148148148149
$ws.Range("K4").Value = 0.4883333333333333
$ws.Range("L4").Value = 0.5263157894736842
$ws.Range("M4").Value = "[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]"
$ws.Range("N4").Value = "[0 1 1 0 0 1 0 1 1 0 0 0 0 1 1 1 0 1 1 0 1 0 1 0]"
$ws.Range("O4").Value = 23

# Row 5
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 BaggingClassifier(estimator=SVC(C=0.25,
                                                 class_weight='balanced',
                                                 random_state=42),
                                   n_estimators=200, random_state=42))])"
$ws.Range("B5").Value = 0.5445487845487846
$ws.Range("C5").Value = "{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': MinMaxScaler(), 'model__n_estimators': 200, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 0.25}"
$ws.Range("D5").Value = 0.6283754338090092
$ws.Range("E5").Value = 0.4328326303326303
$ws.Range("F5").Value = 0.7741935483870968
$ws.Range("G5").Value = 0.8155145468583518
$ws.Range("H5").Value = 0.5731712962962963
$ws.Range("I5").Value = 0.7058823529411765
$ws.Range("J5").Value = 0.5353741496598639
$ws.Range("K5").Value = 0.3930555555555555
$ws.Range("L5").Value = 0.8571428571428571
$ws.Range("M5").Value = "[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]"
$ws.Range("N5").Value = "[0 1 1 1 1 1 0 1 0 0 1 1 1 1 1 0 0 0 1 1 1 1 1 1]"
$ws.Range("O5").Value = 99

# Row 6
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7fa6cd284310>),
                ('model',
                 BaggingClassifier(estimator=SVC(C=0.25,
                                                 class_weight='balanced',
                                                 random_state=42),
                                   n_estimators=200, random_state=42))])"
$ws.Range("B6").Value = 0.6633333333333333
$ws.Range("C6").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f9f104e3e20>, 'scaler': MinMaxScaler(), 'model__n_estimators': 200, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 0.25}"
$ws.Range("D6").Value = 0.7827114645380154
$ws.Range("E6").Value = 0.5938885188885189
$ws.Range("F6").Value = 0.6
$ws.Range("G6").Value = 0.8174656036443545
$ws.Range("H6").Value = 0.6279100529100529
$ws.Range("I6").Value = 0.4736842105263158
$ws.Range("J6").Value = 0.7612179487179488
$ws.Range("K6").Value = 0.5958333333333333
$ws.Range("L6").Value = 0.8181818181818182
$ws.Range("M6").Value = "[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]"
$ws.Range("N6").Value = "[0 1 1 0 0 1 0 1 1 1 1 1 1 0 1 1 0 0 1 0 1 1 1 0]"
$ws.Range("O6").Value = 89
